# Auto-generated edit script: applies 2022-08-30 violent-crime data update.
# For each affected worksheet, update the 2022 (column I) and, in a couple of
# cases, 2015 (column B) totals to their corrected values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 4784
$ws.Range("I3").Value = 4962
$ws.Range("B4").Value = 1658
$ws.Range("I4").Value = 1139
$ws.Range("I5").Value = 455
$ws.Range("I6").Value = 5405
$ws.Range("B7").Value = 23290
$ws.Range("I7").Value = 16745

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I4").Value = 65
$ws.Range("I6").Value = 114
$ws.Range("I7").Value = 534
$ws.Range("I8").Value = 1023
$ws.Range("I9").Value = 76
$ws.Range("I11").Value = 254
$ws.Range("I12").Value = 39
$ws.Range("I14").Value = 94
$ws.Range("I20").Value = 410
$ws.Range("I23").Value = 160
$ws.Range("I27").Value = 151
$ws.Range("I29").Value = 1063
$ws.Range("I31").Value = 164
$ws.Range("I33").Value = 772
$ws.Range("I34").Value = 81
$ws.Range("I35").Value = 22
$ws.Range("I37").Value = 529
$ws.Range("I42").Value = 573
$ws.Range("I43").Value = 135
$ws.Range("I44").Value = 122
$ws.Range("I48").Value = 232
$ws.Range("I49").Value = 125
$ws.Range("I50").Value = 77
$ws.Range("I51").Value = 179
$ws.Range("I52").Value = 359
$ws.Range("I53").Value = 173
$ws.Range("I55").Value = 185
$ws.Range("I56").Value = 17
$ws.Range("I57").Value = 63
$ws.Range("I59").Value = 30
$ws.Range("B63").Value = 363
$ws.Range("I63").Value = 61
$ws.Range("I65").Value = 374
$ws.Range("I66").Value = 48
$ws.Range("I67").Value = 658
$ws.Range("I76").Value = 256
$ws.Range("I79").Value = 469
$ws.Range("I82").Value = 20
$ws.Range("I83").Value = 351
$ws.Range("I85").Value = 756
$ws.Range("I89").Value = 189
$ws.Range("I93").Value = 99
$ws.Range("I94").Value = 162
$ws.Range("I95").Value = 272
$ws.Range("I98").Value = 105
$ws.Range("I99").Value = 314
$ws.Range("B101").Value = 23290
$ws.Range("I101").Value = 16745

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 201
$ws.Range("I4").Value = 41
$ws.Range("I7").Value = 756

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 100
$ws.Range("I7").Value = 359

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I4").Value = 23
$ws.Range("I7").Value = 254

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 324
$ws.Range("I3").Value = 285
$ws.Range("I7").Value = 1023

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 173

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I6").Value = 138
$ws.Range("I7").Value = 534

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I5").Value = 7
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 48
$ws.Range("I6").Value = 63

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 34
$ws.Range("I3").Value = 25
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I6").Value = 143
$ws.Range("I7").Value = 529

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 116
$ws.Range("I7").Value = 314

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 155
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 658

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 121
$ws.Range("I3").Value = 110
$ws.Range("I7").Value = 374

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 128
$ws.Range("I3").Value = 134
$ws.Range("I7").Value = 351

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 95
$ws.Range("I7").Value = 272

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 287
$ws.Range("I6").Value = 243
$ws.Range("I7").Value = 772

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 125

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 314
$ws.Range("I4").Value = 55
$ws.Range("I6").Value = 288
$ws.Range("I7").Value = 1063

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I2").Value = 31
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 232

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I4").Value = 29
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I6").Value = 162
$ws.Range("I7").Value = 573

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 55
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 136
$ws.Range("I3").Value = 151
$ws.Range("I5").Value = 15
$ws.Range("I6").Value = 135
$ws.Range("I7").Value = 469

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 115
$ws.Range("I3").Value = 123
$ws.Range("I6").Value = 129
$ws.Range("I7").Value = 410

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I4").Value = 16
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 30
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 38
$ws.Range("I6").Value = 71
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 20

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("I3").Value = 2
$ws.Range("I7").Value = 17

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 39
